# Remove the "副标题 2" subtitle placeholder shape (id=3, containing the
# "小组成员：吴岳东，刘家溪" text) from slide 1 entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape by its name so the script is resilient to shape
# ordering; falls back to the known index if no name match is found.
$targetIndex = 0
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "副标题 2") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq 0) {
    $targetIndex = 3
}

# This shape is bound to a layout placeholder (subTitle), so the first
# Delete() only clears its content and leaves an empty inherited
# placeholder behind (standard PowerPoint placeholder behaviour). Delete
# it a second time so the shape is fully removed from the slide, matching
# the target OOXML which no longer contains this <p:sp> at all.
$s.Shapes.Item($targetIndex).Delete()
$s.Shapes.Item($targetIndex).Delete()
